# results_Alpha.xlsx — "prove null hypothesis to be rejected using t ttest_ind"
#
# - Subject rows 4-6 (the extra "movehat"/"movehead"/"raisingeyebrows"-only
#   trials) are removed, leaving just the two real trials in rows 2-3.
# - The two remaining trials' Q/R/S (stats) columns get refreshed values
#   from the new t-test (ttest_ind) run.
# - The "Noise" label for each remaining trial is updated to reflect the
#   combined conditions actually tested (eyescrunching+jaw / jaw+raisingeyebrows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 4:6 entirely (shrinks the used range down to A1:U3).
$ws.Rows("4:6").Delete()

# Row 2 — refreshed SNR-before / SNR-after / SNR-LMS stats.
$ws.Range("Q2").Value = -1.712258434503343
$ws.Range("R2").Value = 4.320651653371465
$ws.Range("S2").Value = -2.967762154706217

# Row 3 — refreshed SNR-before / SNR-after / SNR-LMS stats.
$ws.Range("Q3").Value = -1.022131778245286
$ws.Range("R3").Value = 6.42221850955547
$ws.Range("S3").Value = -18.82147087323878

# Relabel the "Noise" condition column for the two surviving trials.
$ws.Range("U2").Value = "eyescrunching+jaw"
$ws.Range("U3").Value = "jaw+raisingeyebrows"
